# Update the "取得日時" (retrieved datetime) column on the active sheet
# ("ランサーズ") for all existing data rows (rows 2-9) to reflect the new
# append timestamp: 2025-10-03 06:33:03 (JST).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2025-10-03 06:33:03"

for ($row = 2; $row -le 9; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
